$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "Name and contact information if an agent is being used."
$ws.Range("B25").Value = "Name and contact information if an agent is being used."
$ws.Range("B33").Value = "Telephone number and email address of the applicant."
$ws.Range("B37").Value = "Name and contact information for the parties making the application."
$ws.Range("B43").Value = "Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation."
$ws.Range("B44").Value = "Details of any conflict of interest that may exist between the applicant and planning authority."
$ws.Range("B47").Value = "Signed and dated verification of the application's accuracy."
$ws.Range("B50").Value = "Details of pre-application advice received from the planning authority"
$ws.Range("B55").Value = "What development, works or change of use is proposed"
$ws.Range("B59").Value = "Where the proposed development will be built."
$ws.Range("B68").Value = "Information to help the planning authority arrange a site visit"
$ws.Range("B74").Value = "Any additional information which will help with the planning application"
